# Update the Bulgaria (BG) row in the "Feeds" sheet: replace the two old
# RSS URLs (url_1/url_2) with a single new URL in url_3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feeds")

$ws.Range("B3").Value = $null
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = "https://www.vesti.bg/rss"

$ws.Activate()
$ws.Range("C11").Select()
